# Technology-Stock-Data.xlsx edit:
#   itemloc sheet: columns D (expiry date) / E (fifo date) / F (max capacity)
#   get rotated so the "max capacity" numbers move from column F into column D,
#   and the two date columns shift right by one (D->E, E->F).
#   Also: itemloc becomes the active sheet/tab (previously salesorder was),
#   with a fresh selection of J8 (no leftover scroll/selection state).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("itemloc")

$lastRow = 63      # header (row 1) + 62 data rows
$colD = 4
$colE = 5
$colF = 6
$scratchCol = 26   # column Z - far outside the used range, used as temp holder

for ($r = 1; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, $colD)
    $eCell = $ws.Cells.Item($r, $colE)
    $fCell = $ws.Cells.Item($r, $colF)
    $scratch = $ws.Cells.Item($r, $scratchCol)

    # stash old D, then rotate F->D, E->F, (old D, from scratch)->E
    $dCell.Copy($scratch) | Out-Null
    $fCell.Copy($dCell) | Out-Null
    $eCell.Copy($fCell) | Out-Null
    $scratch.Copy($eCell) | Out-Null
    $scratch.Clear() | Out-Null
}

# Make itemloc the active sheet/tab (it was "salesorder" before) and set a
# plain single-cell selection, clearing any scrolled topLeftCell state.
$ws.Activate() | Out-Null
$ws.Range("J8").Select() | Out-Null
